# The edit rotates the data of rows 13, 14 and 15 on sheet "Artfynd":
#   new row 13 <- old row 15
#   new row 14 <- old row 13
#   new row 15 <- old row 14
# (row 13's species observation moves down to row 14, row 14's moves down
#  to row 15, and row 15's moves up to row 13 - a cyclic shift).
#
# We overwrite the three rows' cell values directly with what the diff
# shows. Some cells need to end up "present but blank" (they existed as
# empty cells before the edit, or need to become empty placeholders) while
# others need to be removed entirely. Simply assigning Value = "" always
# deletes the cell node, so for cells that must stay present-but-empty we
# clear the value and then touch the Style property (re-assign it to
# itself) which makes the engine keep an (empty) cell node without
# allocating a new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Clear-KeepCell($rng) {
    # Clear a cell's value but keep the <c> element present (empty cell).
    $rng.Value = ""
    $rng.Style = $rng.Style
}

function Clear-RemoveCell($rng) {
    # Clear a cell's value and let the (now valueless) cell be dropped.
    $rng.Value = ""
}

# ---------------- Row 13 (becomes former row 15's data) ----------------
$ws.Range("A13").Value = 112205187
$ws.Range("B13").Value = 90810
$ws.Range("E13").Value = 4363
$ws.Range("F13").Value = "Zontaggsvamp"
$ws.Range("G13").Value = "Hydnellum concrescens"
$ws.Range("H13").Value = "(Pers.) Banker"
$ws.Range("I13").Value = "'1"
$ws.Range("J13").Value = "mycel"
Clear-RemoveCell $ws.Range("L13")
$ws.Range("Q13").Value = 478586
$ws.Range("R13").Value = 6556137
$ws.Range("AI13").Value = "i yta bökad av vildsvin"
Clear-RemoveCell $ws.Range("AJ13")
Clear-RemoveCell $ws.Range("AK13")
Clear-RemoveCell $ws.Range("AM13")
Clear-RemoveCell $ws.Range("AO13")

# ---------------- Row 14 (becomes former row 13's data) ----------------
$ws.Range("A14").Value = 112231588
$ws.Range("B14").Value = 93334
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 2818
$ws.Range("F14").Value = "Stubbspretmossa"
$ws.Range("G14").Value = "Herzogiella seligeri"
$ws.Range("H14").Value = "(Brid.) Z.Iwats."
Clear-KeepCell $ws.Range("I14")
Clear-KeepCell $ws.Range("J14")
Clear-RemoveCell $ws.Range("M14")
$ws.Range("Q14").Value = 478719
$ws.Range("R14").Value = 6556487
Clear-KeepCell $ws.Range("AF14")
$ws.Range("AH14").Value = "Blåbärsbarrskog"
$ws.Range("AJ14").Value = "tall"
$ws.Range("AK14").Value = "Pinus sylvestris"
$ws.Range("AM14").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO14").Value = "Horizontal, dead with ground contact # murken grov låga # Pinus sylvestris"

# ---------------- Row 15 (becomes former row 14's data) ----------------
$ws.Range("A15").Value = 112231491
$ws.Range("B15").Value = 56575
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 103021
$ws.Range("F15").Value = "Talltita"
$ws.Range("G15").Value = "Poecile montanus"
$ws.Range("H15").Value = "(Conrad von Baldenstein, 1827)"
Clear-RemoveCell $ws.Range("J15")
Clear-KeepCell $ws.Range("L15")
$ws.Range("M15").Value = "lockläte, övriga läten"
$ws.Range("Q15").Value = 478579
$ws.Range("R15").Value = 6556322
Clear-RemoveCell $ws.Range("AF15")
Clear-RemoveCell $ws.Range("AH15")
Clear-RemoveCell $ws.Range("AI15")

$wb.Save()
